# km-changes.docx: strip stray spell/grammar-check markers that had split
# the "ALTER TABLE ..." sentences into multiple runs, collapsing each back
# into a single run of text per sentence. Then append the new migration
# statements for sizes/product_attributes/pages/faqs/campaigns tables.

$d = $word.ActiveDocument

function Replace-Plain($text) {
    $d.Content.Find.Execute($text, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $text, 2) | Out-Null
}

# Paragraphs 1, 2, 4-10: each is one sentence split across several runs by
# w:proofErr spell/grammar markers -- re-"typing" the same text over the
# whole match merges it back into a single clean run.
Replace-Plain "ALTER TABLE products ADD name_km INT NULL DEFAULT NULL AFTER name;"
Replace-Plain "ALTER TABLE products ADD summary_km VARCHAR(255) NULL DEFAULT NULL AFTER summary;"

# Paragraph 3 holds two sentences separated by two manual line breaks
# (w:br). Replace each sentence's text independently so the br runs in
# between are left untouched.
Replace-Plain "ALTER TABLE products ADD description_km LONGTEXT NULL DEFAULT NULL AFTER description;"
Replace-Plain "ALTER TABLE categories ADD name_km VARCHAR(255) NULL DEFAULT NULL AFTER name;"

Replace-Plain "ALTER TABLE sub_categories ADD name_km VARCHAR(255) NULL DEFAULT NULL AFTER name;"
Replace-Plain "ALTER TABLE units ADD name_km VARCHAR(255) NULL DEFAULT NULL AFTER name;"
Replace-Plain "ALTER TABLE delivery_options ADD title_km VARCHAR(255) NULL DEFAULT NULL AFTER title;"
Replace-Plain "ALTER TABLE delivery_options ADD sub_title_km VARCHAR(255) NULL DEFAULT NULL AFTER sub_title;"
Replace-Plain "ALTER TABLE brands ADD name_km VARCHAR(255) NULL DEFAULT NULL AFTER name;"
Replace-Plain "ALTER TABLE brands ADD description_km LONGTEXT NULL DEFAULT NULL AFTER description;"
Replace-Plain "ALTER TABLE colors ADD name_km VARCHAR(255) NULL DEFAULT NULL AFTER name;"

# New trailing paragraphs (one ALTER statement per paragraph) appended
# after the "colors" line.
$newStatements = @(
    "ALTER TABLE sizes ADD name_km VARCHAR(255) NULL DEFAULT NULL AFTER name;",
    "ALTER TABLE product_attributes ADD title_km VARCHAR(255) NULL DEFAULT NULL AFTER title;",
    "ALTER TABLE product_attributes ADD terms_km LONGTEXT NULL DEFAULT NULL AFTER terms;",
    "ALTER TABLE pages ADD title_km VARCHAR(255) NULL DEFAULT NULL AFTER title;",
    "ALTER TABLE pages ADD meta_tags_km TEXT NULL DEFAULT NULL AFTER meta_tags;",
    "ALTER TABLE pages ADD content_km LONGTEXT NULL DEFAULT NULL AFTER content;",
    "ALTER TABLE pages ADD meta_description_km TEXT NULL DEFAULT NULL AFTER meta_description;",
    "ALTER TABLE faqs ADD title_km VARCHAR(255) NULL DEFAULT NULL AFTER title;",
    "ALTER TABLE faqs ADD description_km LONGTEXT NULL DEFAULT NULL AFTER description;",
    "ALTER TABLE campaigns ADD title_km VARCHAR(255) NULL DEFAULT NULL AFTER title;",
    "ALTER TABLE campaigns CHANGE subtitle subtitle TEXT CHARACTER SET utf8mb4 COLLATE utf8mb4_unicode_ci NULL DEFAULT NULL;",
    "ALTER TABLE campaigns ADD subtitle_km TEXT NULL DEFAULT NULL AFTER subtitle;"
)

foreach ($stmt in $newStatements) {
    $tail = $d.Paragraphs.Item($d.Paragraphs.Count).Range
    $tail.Collapse(0)
    $tail.InsertParagraphAfter()
    $d.Paragraphs.Item($d.Paragraphs.Count).Range.Text = $stmt
}
